$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Assign values in the exact order needed so that the shared-strings table
# ends up built in the same first-use order as the target workbook:
#   turnos, decorre, preferencia, Webdev, Num Pessoas,
#   Turno1, Turno2, Turno3, Turno4, Turno5, [1,1,1,1,1,1,1,1,1,1]
$ws.Range("A1").Value = "turnos"
$ws.Range("B1").Value = "decorre"
$ws.Range("D1").Value = "preferencia"
$ws.Range("D2").Value = "Webdev"
$ws.Range("C1").Value = "Num Pessoas"
$ws.Range("A2").Value = "Turno1"
$ws.Range("A3").Value = "Turno2"
$ws.Range("A4").Value = "Turno3"
$ws.Range("A5").Value = "Turno4"
$ws.Range("A6").Value = "Turno5"
$ws.Range("B2").Value = "[1,1,1,1,1,1,1,1,1,1]"

# Remaining cells reuse already-registered shared strings / plain numbers
$ws.Range("B3").Value = "[1,1,1,1,1,1,1,1,1,1]"
$ws.Range("B4").Value = "[1,1,1,1,1,1,1,1,1,1]"
$ws.Range("B5").Value = "[1,1,1,1,1,1,1,1,1,1]"
$ws.Range("B6").Value = "[1,1,1,1,1,1,1,1,1,1]"

$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 3

$ws.Range("D3").Value = "Webdev"
$ws.Range("D4").Value = "Webdev"
$ws.Range("D5").Value = "Webdev"
$ws.Range("D6").Value = "Webdev"

# Adjust column B width to fit the new, wider content and select B2 as active cell
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Range("B2").Select() | Out-Null

$wb.Save()
